# The deck's design was switched from the "Integral" (Red Violet) theme
# back to the default "Office Theme" palette. Re-point every themed
# colour slot on the slide master to the stock Office colour scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeRGB {
    param($index, $red, $green, $blue)
    $colorScheme.Item($index).RGB = $red + ($green * 256) + ($blue * 65536)
}

# Office Theme colour scheme, in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
Set-ThemeRGB 1  0x00 0x00 0x00   # dk1
Set-ThemeRGB 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeRGB 3  0x44 0x54 0x6A   # dk2
Set-ThemeRGB 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRGB 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRGB 6  0xED 0x7D 0x31   # accent2
Set-ThemeRGB 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRGB 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRGB 9  0x44 0x72 0xC4   # accent5
Set-ThemeRGB 10 0x70 0xAD 0x47   # accent6
Set-ThemeRGB 11 0x05 0x63 0xC1   # hlink
Set-ThemeRGB 12 0x95 0x4F 0x72   # folHlink
